# Daily update at 8 AM UTC
# Append the next day's row of win totals to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 46017
$ws.Cells.Item($newRow, 2).Value = 148
$ws.Cells.Item($newRow, 3).Value = 160
$ws.Cells.Item($newRow, 4).Value = 149

# Match the date-style formatting used by the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
